# Applies the "Updated cryptos list" data refresh:
#   - per-row Price (D) / Volume(1h) (E) text updates
#   - rows 41/42 swap (Stacks <-> Stellar) with their new data
#
# Every value is written with a leading apostrophe (Excel's literal
# "treat as text" prefix) so look-alike numbers (e.g. "357.46" or
# "51.928.58") are stored as text, matching the sheet's inlineStr/shared-
# string cells instead of being parsed into numeric cells. The style is
# then reset to "Normal" so the quote-prefix flag does not leave a stray
# cell-format change behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = 'Normal'
}

# Row 2
Set-TextValue 'D2' '51.928.58'
Set-TextValue 'E2' '  -0.02%  '

# Row 3
Set-TextValue 'D3' '2.777.67'
Set-TextValue 'E3' '  -1.60%  '

# Row 4
Set-TextValue 'E4' '  -0.01%  '

# Row 5
Set-TextValue 'D5' '357.46'
Set-TextValue 'E5' '  +0.59%  '

# Row 6
Set-TextValue 'D6' '108.96'
Set-TextValue 'E6' '  -4.04%  '

# Row 7
Set-TextValue 'D7' '0.567'
Set-TextValue 'E7' '  +2.60%  '

# Row 9
Set-TextValue 'E9' '  -1.47%  '

# Row 10
Set-TextValue 'D10' '40.08'
Set-TextValue 'E10' '  -4.00%  '

# Row 11
Set-TextValue 'D11' '0.0851'
Set-TextValue 'E11' '  -0.21%  '

# Row 12
Set-TextValue 'E12' '  +0.72%  '

# Row 13
Set-TextValue 'D13' '19.37'
Set-TextValue 'E13' '  -3.31%  '

# Row 14
Set-TextValue 'E14' '  -1.44%  '

# Row 15
Set-TextValue 'D15' '3.215.73'
Set-TextValue 'E15' '  -0.94%  '

# Row 16
Set-TextValue 'D16' '2.804.66'
Set-TextValue 'E16' '  -1.09%  '

# Row 17
Set-TextValue 'D17' '0.929'
Set-TextValue 'E17' '  +3.39%  '

# Row 18
Set-TextValue 'D18' '51.873.68'
Set-TextValue 'E18' '  +0.02%  '

# Row 19
Set-TextValue 'D19' '7.40'
Set-TextValue 'E19' '  -0.11%  '

# Row 20
Set-TextValue 'D20' '3.13'
Set-TextValue 'E20' '  -0.87%  '

# Row 21
Set-TextValue 'D21' '13.03'
Set-TextValue 'E21' '  -4.46%  '

# Row 22
Set-TextValue 'E22' '  -2.05%  '

# Row 23
Set-TextValue 'D23' '274.15'
Set-TextValue 'E23' '  +1.53%  '

# Row 24
Set-TextValue 'E24' '  +0.02%  '

# Row 25
Set-TextValue 'D25' '2.74'
Set-TextValue 'E25' '  -2.01%  '

# Row 26
Set-TextValue 'D26' '26.53'
Set-TextValue 'E26' '  -0.98%  '

# Row 27
Set-TextValue 'E27' '  -0.13%  '

# Row 28
Set-TextValue 'E28' '  -1.79%  '

# Row 29
Set-TextValue 'E29' '  -1.24%  '

# Row 30
Set-TextValue 'E30' '  +1.73%  '

# Row 31
Set-TextValue 'D31' '0.0465'
Set-TextValue 'E31' '  +1.21%  '

# Row 32
Set-TextValue 'D32' '51.69'
Set-TextValue 'E32' '  +1.50%  '

# Row 33
Set-TextValue 'D33' '33.76'
Set-TextValue 'E33' '  -0.24%  '

# Row 34
Set-TextValue 'D34' '5.70'
Set-TextValue 'E34' '  -2.28%  '

# Row 35
Set-TextValue 'E35' '  +1.44%  '

# Row 36
Set-TextValue 'D36' '5.28'
Set-TextValue 'E36' '  +7.39%  '

# Row 37
Set-TextValue 'E37' '  +0.03%  '

# Row 38
Set-TextValue 'E38' '  +0.40%  '

# Row 39
Set-TextValue 'D39' '18.05'
Set-TextValue 'E39' '  -1.84%  '

# Row 40
Set-TextValue 'E40' '  -4.74%  '

# Row 41
Set-TextValue 'B41' 'Stellar'
Set-TextValue 'C41' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D41' '0.114'
Set-TextValue 'E41' '  -0.68%  '

# Row 42
Set-TextValue 'B42' 'Stacks'
Set-TextValue 'C42' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D42' '2.51'
Set-TextValue 'E42' '  -1.86%  '

# Row 43
Set-TextValue 'E43' '  -2.92%  '

# Row 44
Set-TextValue 'D44' '121.33'
Set-TextValue 'E44' '  -5.76%  '

# Row 45
Set-TextValue 'D45' '21.91'
Set-TextValue 'E45' '  -6.91%  '

# Row 46
Set-TextValue 'D46' '2.065.18'
Set-TextValue 'E46' '  -0.71%  '

# Row 47
Set-TextValue 'D47' '3.24'
Set-TextValue 'E47' '  -3.57%  '

# Row 48
Set-TextValue 'D48' '2.26'
Set-TextValue 'E48' '  -2.06%  '

# Row 49
Set-TextValue 'D49' '5.69'
Set-TextValue 'E49' '  +0.03%  '

# Row 50
Set-TextValue 'D50' '0.921'

# Row 51
Set-TextValue 'D51' '8.95'
Set-TextValue 'E51' '  +0.23%  '
